# Revised ACF popsize calculation again
#
# Insert a new "program_prop_population_screened" parameter row into the
# "constants" sheet. It belongs right after the existing
# "program_ratio_case_detection_with_raised_awareness" row (old row 58),
# continuing that program_* block, and pushes every row below it down by
# one - exactly what Excel's own Insert Row does.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new row at position 58; everything from the old row 58 onward
# (including row formatting, dataValidation ranges, etc.) shifts down to
# 59:170, and the new blank row inherits the look of the row above it.
$ws.Rows.Item(58).Insert() | Out-Null

# Populate the freshly inserted row with the new parameter name + its
# default value.
$ws.Range("A58").Value = "program_prop_population_screened"
$ws.Range("B58").Value = 0.8

# Leave the workbook scrolled/selected where the author ended up after
# making the edit.
$ws.Activate() | Out-Null
$ws.Range("A50").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 50
$ws.Range("A56").Select() | Out-Null
